$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 20:23"

# Row 4
$ws.Range("B4").Value = 5579686
$ws.Range("C4").Value = 13054
$ws.Range("D4").Value = 2928179
$ws.Range("E4").Value = 2478249
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 130
$ws.Range("H4").Value = 173258

# Row 6
$ws.Range("B6").Value = 2694614
$ws.Range("C6").Value = 47298
$ws.Range("D6").Value = 1968175
$ws.Range("E6").Value = 674632
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 762
$ws.Range("H6").Value = 51807

# Row 21
$ws.Range("B21").Value = 250542
$ws.Range("C21").Value = 1233
$ws.Range("D21").Value = 231971
$ws.Range("E21").Value = 12575
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 5996

# Row 24
$ws.Range("B24").Value = 180133
$ws.Range("C24").Value = 3202
$ws.Range("D24").Value = 128945
$ws.Range("E24").Value = 45234
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 5954

# Row 30
$ws.Range("B30").Value = 101751
$ws.Range("C30").Value = 209
$ws.Range("D30").Value = 87120
$ws.Range("E30").Value = 8548
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 6083

# Row 64
$ws.Range("A64").Value = "Etiopia"
$ws.Range("B64").Value = 31336
$ws.Range("C64").Value = 1460
$ws.Range("D64").Value = 12524
$ws.Range("E64").Value = 18268
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 16
$ws.Range("H64").Value = 544

# Row 65
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 30377
$ws.Range("C65").Value = 194
$ws.Range("D65").Value = 21220
$ws.Range("E65").Value = 8249
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 12
$ws.Range("H65").Value = 908

# Row 66
$ws.Range("A66").Value = "Kenia"
$ws.Range("B66").Value = 30365
$ws.Range("C66").Value = 245
$ws.Range("D66").Value = 17160
$ws.Range("E66").Value = 12723
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 482

# Row 69
$ws.Range("B69").Value = 27313
$ws.Range("C69").Value = 56
$ws.Range("D69").Value = 23364
$ws.Range("E69").Value = 2175
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 1774

# Row 77
$ws.Range("B77").Value = 16844
$ws.Range("C77").Value = 310
$ws.Range("D77").Value = 9906
$ws.Range("E77").Value = 6825
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 113

# Row 87
$ws.Range("A87").Value = "Zambia"
$ws.Range("B87").Value = 9839
$ws.Range("C87").Value = 496
$ws.Range("D87").Value = 8575
$ws.Range("E87").Value = 1000
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 264

# Row 88
$ws.Range("A88").Value = "Paraguay"
$ws.Range("B88").Value = 9791
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 6034
$ws.Range("E88").Value = 3619
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 138

# Row 89
$ws.Range("A89").Value = "Consejo Danes para los Refugiados"
$ws.Range("B89").Value = 9706
$ws.Range("C89").Value = 30
$ws.Range("D89").Value = 8705
$ws.Range("E89").Value = 758
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 243

# Row 99
$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 7499
$ws.Range("C99").Value = 119
$ws.Range("D99").Value = 3816
$ws.Range("E99").Value = 3453
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 230

# Row 100
$ws.Range("A100").Value = "Luxemburgo"
$ws.Range("B100").Value = 7458
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 6500
$ws.Range("E100").Value = 835
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 123

# Row 123
$ws.Range("A123").Value = "Mozambique"
$ws.Range("B123").Value = 2914
$ws.Range("C123").Value = 59
$ws.Range("D123").Value = 1196
$ws.Range("E123").Value = 1699
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 19

# Row 124
$ws.Range("A124").Value = "Eslovaquia"
$ws.Range("B124").Value = 2907
$ws.Range("C124").Value = 5
$ws.Range("D124").Value = 1969
$ws.Range("E124").Value = 907
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 31

# Row 125
$ws.Range("A125").Value = "Sri Lanka"
$ws.Range("B125").Value = 2900
$ws.Range("C125").Value = 7
$ws.Range("D125").Value = 2676
$ws.Range("E125").Value = 213
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 11

# Row 162
$ws.Range("B162").Value = 880
$ws.Range("C162").Value = 25
$ws.Range("D162").Value = 657
$ws.Range("E162").Value = 218
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 5

# Row 167
$ws.Range("A167").Value = "Guadalupe"
$ws.Range("B167").Value = 510
$ws.Range("C167").Value = 64
$ws.Range("D167").Value = 289
$ws.Range("E167").Value = 207
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 14

# Row 168
$ws.Range("A168").Value = "Tanzania"
$ws.Range("B168").Value = 509
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 183
$ws.Range("E168").Value = 305
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 21

# Row 169
$ws.Range("A169").Value = "Taiwan"
$ws.Range("B169").Value = 485
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 450
$ws.Range("E169").Value = 28
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 7

# Row 170
$ws.Range("A170").Value = "Belice"
$ws.Range("B170").Value = 452
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 35
$ws.Range("E170").Value = 414
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 3

# Row 188
$ws.Range("B188").Value = 152
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 122
$ws.Range("E188").Value = 23
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 7

# Row 213
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
